$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.057.56'
$ws.Range('E2').Value = '  -3.31%  '

$ws.Range('D3').Value = '3.076.57'
$ws.Range('E3').Value = '  -1.79%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').Value = '3.069.62'
$ws.Range('E8').Value = '  -1.40%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.158'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.56%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.52'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '34.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.56%  '

$ws.Range('E14').Value = '  -2.08%  '

$ws.Range('D15').Value = '3.572.60'
$ws.Range('E15').Value = '  -1.65%  '

$ws.Range('D16').Value = '63.118.82'
$ws.Range('E16').Value = '  -3.22%  '

$ws.Range('E17').Value = '  -1.24%  '

$ws.Range('D18').Value = '3.078.33'
$ws.Range('E18').Value = '  -1.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '500.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.55%  '

$ws.Range('E22').Value = '  -0.49%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.70%  '

$ws.Range('E25').Value = '  -3.74%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.36%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.62%  '

$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('E30').Value = '  -9.63%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.20%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.42%  '

$ws.Range('E33').Value = '  -6.11%  '

$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.63%  '

$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '523.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.79%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.92'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.69%  '

$ws.Range('E37').Value = '  -6.47%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0400'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.34%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0790'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.81%  '

$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.044.85'
$ws.Range('E40').Value = '  -0.97%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.118'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.54%  '

$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.52%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.00%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.256'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.75%  '

$ws.Range('E46').Value = '  -7.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.76%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.16%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.107'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.51%  '

$ws.Range('B50').Value = 'CoreDAO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.45'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +66.70%  '

$ws.Range('D51').Value = '0.0₃0502'
$ws.Range('E51').Value = '  -5.64%  '
